$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

$indicators = @("INDICATOR_50","INDICATOR_51","INDICATOR_52","INDICATOR_53","INDICATOR_54","INDICATOR_200","INDICATOR_201","INDICATOR_202","INDICATOR_203","INDICATOR_204","INDICATOR_205")

$row = 42
foreach ($ind in $indicators) {
    $ws.Cells.Item(41, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $ws.Cells.Item(41, 2).Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)

    $ws.Cells.Item(41, 3).Copy()
    $ws.Cells.Item($row, 3).PasteSpecial(-4122)

    $ws.Cells.Item(41, 5).Copy()
    $ws.Cells.Item($row, 5).PasteSpecial(-4122)

    $ws.Cells.Item(41, 6).Copy()
    $ws.Cells.Item($row, 6).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = "CREATE/MODIFY"
    $ws.Cells.Item($row, 2).Value = "LIB_EWS_IT"
    $ws.Cells.Item($row, 3).Value = $ind
    $ws.Cells.Item($row, 5).Value = "String"
    $ws.Cells.Item($row, 6).Value = "String"
    $row = $row + 1
}

$excel.CutCopyMode = $false

$ws.Range("A52").Select()

$labels = $wb.Worksheets.Item("Formula Libraries Labels")
$labels.Select()
